$d = $word.ActiveDocument

# The document's header/footer each carry two Pearson/BTec logo pictures
# (primary + first-page variants). The picture's cosmetic internal
# "name" label (visible on the <wp:docPr>/<pic:cNvPr> pair in the
# underlying XML) needs swapping between the two images:
#   - BTec_Logo-Orange (inline picture in the headers):  image1.jpg -> image2.jpg
#   - PearsonLogo.png  (inline picture in the footers):   image2.png -> image1.png
#
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

$sec = $d.Sections(1)

foreach ($hfIndex in 1, 2) {
    $hdr = $sec.Headers($hfIndex)
    if ($hdr.Range.InlineShapes.Count -ge 1) {
        $shp = $hdr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }

    $ftr = $sec.Footers($hfIndex)
    if ($ftr.Range.InlineShapes.Count -ge 1) {
        $shp = $ftr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -like "*PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}
